$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove existing hyperlinks (also clears all hyperlink relationships in one shot)
$ws.Hyperlinks.Delete()

# Delete row 18 (Joe Wieskamp no longer on the roster)
$ws.Rows(18).Delete()

# Force the "Exp" column to keep storing numeric-looking values ("1","5",...) as text
$ws.Range("I2:I17").NumberFormat = "@"

# Update roster data rows 2-17 with the refreshed 2023 roster information
# Row 2: Scottie Barnes
$ws.Cells.Item(2, 2).Value2 = 10
$ws.Cells.Item(2, 3).Value2 = "Scottie Barnes"
$ws.Cells.Item(2, 4).Value2 = "PF"
$ws.Cells.Item(2, 5).Value2 = "6-9"
$ws.Cells.Item(2, 6).Value2 = 227
$ws.Cells.Item(2, 7).Value2 = "August 1, 2001"
$ws.Cells.Item(2, 8).Value2 = "us"
$ws.Cells.Item(2, 9).Value2 = "1"
$ws.Cells.Item(2, 10).Value2 = "Florida State"

# Row 3: Chris Boucher
$ws.Cells.Item(3, 2).Value2 = 25
$ws.Cells.Item(3, 3).Value2 = "Chris Boucher"
$ws.Cells.Item(3, 4).Value2 = "PF"
$ws.Cells.Item(3, 5).Value2 = "6-9"
$ws.Cells.Item(3, 6).Value2 = 200
$ws.Cells.Item(3, 7).Value2 = "January 11, 1993"
$ws.Cells.Item(3, 8).Value2 = "lc"
$ws.Cells.Item(3, 9).Value2 = "5"
$ws.Cells.Item(3, 10).Value2 = "Oregon"

# Row 4: Gary Trent Jr.
$ws.Cells.Item(4, 2).Value2 = 33
$ws.Cells.Item(4, 3).Value2 = "Gary Trent Jr."
$ws.Cells.Item(4, 4).Value2 = "SG"
$ws.Cells.Item(4, 5).Value2 = "6-5"
$ws.Cells.Item(4, 6).Value2 = 209
$ws.Cells.Item(4, 7).Value2 = "January 18, 1999"
$ws.Cells.Item(4, 8).Value2 = "us"
$ws.Cells.Item(4, 9).Value2 = "4"
$ws.Cells.Item(4, 10).Value2 = "Duke"

# Row 5: Fred VanVleet
$ws.Cells.Item(5, 2).Value2 = 23
$ws.Cells.Item(5, 3).Value2 = "Fred VanVleet"
$ws.Cells.Item(5, 4).Value2 = "PG"
$ws.Cells.Item(5, 5).Value2 = "6-1"
$ws.Cells.Item(5, 6).Value2 = 197
$ws.Cells.Item(5, 7).Value2 = "February 25, 1994"
$ws.Cells.Item(5, 8).Value2 = "us"
$ws.Cells.Item(5, 9).Value2 = "6"
$ws.Cells.Item(5, 10).Value2 = "Wichita State"

# Row 6: Pascal Siakam
$ws.Cells.Item(6, 2).Value2 = 43
$ws.Cells.Item(6, 3).Value2 = "Pascal Siakam"
$ws.Cells.Item(6, 4).Value2 = "C"
$ws.Cells.Item(6, 5).Value2 = "6-9"
$ws.Cells.Item(6, 6).Value2 = 230
$ws.Cells.Item(6, 7).Value2 = "April 2, 1994"
$ws.Cells.Item(6, 8).Value2 = "cm"
$ws.Cells.Item(6, 9).Value2 = "6"
$ws.Cells.Item(6, 10).Value2 = "New Mexico State"

# Row 7: Thaddeus Young
$ws.Cells.Item(7, 2).Value2 = 21
$ws.Cells.Item(7, 3).Value2 = "Thaddeus Young"
$ws.Cells.Item(7, 4).Value2 = "PF"
$ws.Cells.Item(7, 5).Value2 = "6-8"
$ws.Cells.Item(7, 6).Value2 = 235
$ws.Cells.Item(7, 7).Value2 = "June 21, 1988"
$ws.Cells.Item(7, 8).Value2 = "us"
$ws.Cells.Item(7, 9).Value2 = "15"
$ws.Cells.Item(7, 10).Value2 = "Georgia Tech"

# Row 8: OG Anunoby
$ws.Cells.Item(8, 2).Value2 = 3
$ws.Cells.Item(8, 3).Value2 = "OG Anunoby"
$ws.Cells.Item(8, 4).Value2 = "SF"
$ws.Cells.Item(8, 5).Value2 = "6-7"
$ws.Cells.Item(8, 6).Value2 = 232
$ws.Cells.Item(8, 7).Value2 = "July 17, 1997"
$ws.Cells.Item(8, 8).Value2 = "gb"
$ws.Cells.Item(8, 9).Value2 = "5"
$ws.Cells.Item(8, 10).Value2 = "Indiana"

# Row 9: Christian Koloko
$ws.Cells.Item(9, 2).Value2 = 35
$ws.Cells.Item(9, 3).Value2 = "Christian Koloko"
$ws.Cells.Item(9, 4).Value2 = "C"
$ws.Cells.Item(9, 5).Value2 = "7-1"
$ws.Cells.Item(9, 6).Value2 = 230
$ws.Cells.Item(9, 7).Value2 = "June 20, 2000"
$ws.Cells.Item(9, 8).Value2 = "cm"
$ws.Cells.Item(9, 9).Value2 = "R"
$ws.Cells.Item(9, 10).Value2 = "Arizona"

# Row 10: Malachi Flynn
$ws.Cells.Item(10, 2).Value2 = 22
$ws.Cells.Item(10, 3).Value2 = "Malachi Flynn"
$ws.Cells.Item(10, 4).Value2 = "PG"
$ws.Cells.Item(10, 5).Value2 = "6-1"
$ws.Cells.Item(10, 6).Value2 = 175
$ws.Cells.Item(10, 7).Value2 = "May 10, 1998"
$ws.Cells.Item(10, 8).Value2 = "us"
$ws.Cells.Item(10, 9).Value2 = "2"
$ws.Cells.Item(10, 10).Value2 = "Washington State, San Diego State"

# Row 11: Juancho Hernangómez
$ws.Cells.Item(11, 2).Value2 = 41
$ws.Cells.Item(11, 3).Value2 = "Juancho Hernangómez"
$ws.Cells.Item(11, 4).Value2 = "PF"
$ws.Cells.Item(11, 5).Value2 = "6-9"
$ws.Cells.Item(11, 6).Value2 = 214
$ws.Cells.Item(11, 7).Value2 = "September 28, 1995"
$ws.Cells.Item(11, 8).Value2 = "es"
$ws.Cells.Item(11, 9).Value2 = "6"
$ws.Cells.Item(11, 10).Value2 = ""

# Row 12: Precious Achiuwa
$ws.Cells.Item(12, 2).Value2 = 5
$ws.Cells.Item(12, 3).Value2 = "Precious Achiuwa"
$ws.Cells.Item(12, 4).Value2 = "C"
$ws.Cells.Item(12, 5).Value2 = "6-8"
$ws.Cells.Item(12, 6).Value2 = 225
$ws.Cells.Item(12, 7).Value2 = "September 19, 1999"
$ws.Cells.Item(12, 8).Value2 = "ng"
$ws.Cells.Item(12, 9).Value2 = "2"
$ws.Cells.Item(12, 10).Value2 = "Memphis"

# Row 13: Dalano Banton
$ws.Cells.Item(13, 2).Value2 = 45
$ws.Cells.Item(13, 3).Value2 = "Dalano Banton"
$ws.Cells.Item(13, 4).Value2 = "PG"
$ws.Cells.Item(13, 5).Value2 = "6-9"
$ws.Cells.Item(13, 6).Value2 = 204
$ws.Cells.Item(13, 7).Value2 = "November 7, 1999"
$ws.Cells.Item(13, 8).Value2 = "ca"
$ws.Cells.Item(13, 9).Value2 = "1"
$ws.Cells.Item(13, 10).Value2 = "Nebraska"

# Row 14: Jeff Dowtin (TW)
$ws.Cells.Item(14, 2).Value2 = 20
$ws.Cells.Item(14, 3).Value2 = "Jeff Dowtin (TW)"
$ws.Cells.Item(14, 4).Value2 = "PG"
$ws.Cells.Item(14, 5).Value2 = "6-3"
$ws.Cells.Item(14, 6).Value2 = 185
$ws.Cells.Item(14, 7).Value2 = "May 10, 1997"
$ws.Cells.Item(14, 8).Value2 = "us"
$ws.Cells.Item(14, 9).Value2 = "1"
$ws.Cells.Item(14, 10).Value2 = "Rhode Island"

# Row 15: Otto Porter Jr.
$ws.Cells.Item(15, 2).Value2 = 32
$ws.Cells.Item(15, 3).Value2 = "Otto Porter Jr."
$ws.Cells.Item(15, 4).Value2 = "SF"
$ws.Cells.Item(15, 5).Value2 = "6-8"
$ws.Cells.Item(15, 6).Value2 = 198
$ws.Cells.Item(15, 7).Value2 = "June 3, 1993"
$ws.Cells.Item(15, 8).Value2 = "us"
$ws.Cells.Item(15, 9).Value2 = "9"
$ws.Cells.Item(15, 10).Value2 = "Georgetown"

# Row 16: Ron Harper Jr. (TW)
$ws.Cells.Item(16, 2).Value2 = 8
$ws.Cells.Item(16, 3).Value2 = "Ron Harper Jr. (TW)"
$ws.Cells.Item(16, 4).Value2 = "SF"
$ws.Cells.Item(16, 5).Value2 = "6-6"
$ws.Cells.Item(16, 6).Value2 = 245
$ws.Cells.Item(16, 7).Value2 = "April 12, 2000"
$ws.Cells.Item(16, 8).Value2 = "us"
$ws.Cells.Item(16, 9).Value2 = "R"
$ws.Cells.Item(16, 10).Value2 = "Rutgers University"

# Row 17: Jakob Poeltl
$ws.Cells.Item(17, 2).Value2 = 19
$ws.Cells.Item(17, 3).Value2 = "Jakob Poeltl"
$ws.Cells.Item(17, 4).Value2 = "C"
$ws.Cells.Item(17, 5).Value2 = "7-1"
$ws.Cells.Item(17, 6).Value2 = 245
$ws.Cells.Item(17, 7).Value2 = "October 15, 1995"
$ws.Cells.Item(17, 8).Value2 = "at"
$ws.Cells.Item(17, 9).Value2 = "6"
$ws.Cells.Item(17, 10).Value2 = "Utah"

# Re-create hyperlinks for column K (bbref url) in row order so relationship ids line up
$ws.Hyperlinks.Add($ws.Range("K2"), "https://www.basketball-reference.com/players/b/barnesc01.html") | Out-Null
$ws.Hyperlinks.Add($ws.Range("K3"), "https://www.basketball-reference.com/players/b/bouchch01.html") | Out-Null
$ws.Hyperlinks.Add($ws.Range("K4"), "https://www.basketball-reference.com/players/t/trentga02.html") | Out-Null
$ws.Hyperlinks.Add($ws.Range("K5"), "https://www.basketball-reference.com/players/v/vanvlfr01.html") | Out-Null
$ws.Hyperlinks.Add($ws.Range("K6"), "https://www.basketball-reference.com/players/s/siakapa01.html") | Out-Null
$ws.Hyperlinks.Add($ws.Range("K7"), "https://www.basketball-reference.com/players/y/youngth01.html") | Out-Null
$ws.Hyperlinks.Add($ws.Range("K8"), "https://www.basketball-reference.com/players/a/anunoog01.html") | Out-Null
$ws.Hyperlinks.Add($ws.Range("K9"), "https://www.basketball-reference.com/players/k/kolokch01.html") | Out-Null
$ws.Hyperlinks.Add($ws.Range("K10"), "https://www.basketball-reference.com/players/f/flynnma01.html") | Out-Null
$ws.Hyperlinks.Add($ws.Range("K11"), "https://www.basketball-reference.com/players/h/hernaju01.html") | Out-Null
$ws.Hyperlinks.Add($ws.Range("K12"), "https://www.basketball-reference.com/players/a/achiupr01.html") | Out-Null
$ws.Hyperlinks.Add($ws.Range("K13"), "https://www.basketball-reference.com/players/b/bantoda01.html") | Out-Null
$ws.Hyperlinks.Add($ws.Range("K14"), "https://www.basketball-reference.com/players/d/dowtije01.html") | Out-Null
$ws.Hyperlinks.Add($ws.Range("K15"), "https://www.basketball-reference.com/players/p/porteot01.html") | Out-Null
$ws.Hyperlinks.Add($ws.Range("K16"), "https://www.basketball-reference.com/players/h/harpero02.html") | Out-Null
$ws.Hyperlinks.Add($ws.Range("K17"), "https://www.basketball-reference.com/players/p/poeltja01.html") | Out-Null

# Re-apply the Hyperlink cell style (Add() creates a fresh duplicate style; restore the original one)
$ws.Range("K2:K17").Style = "Hyperlink"
